$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Forces the cell to keep its value as literal text even when the
    # string looks like a number (e.g. "244.64" or "17.50"), without
    # leaving a lingering "Text" number-format style on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "36.293.08"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.040.32"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "244.64"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - Solana
Set-TextValue "D8" "54.57"
$ws.Range("E8").Value = "  +1.81%  "

# Row 9 - OKB
Set-TextValue "D9" "59.58"
$ws.Range("E9").Value = "  +1.29%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.363"
$ws.Range("E10").Value = "  -0.54%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0738"
$ws.Range("E11").Value = "  -2.89%  "

# Row 13 - Polygon
Set-TextValue "D13" "0.898"
$ws.Range("E13").Value = "  +1.42%  "

# Row 14 - Chainlink
Set-TextValue "D14" "14.27"
$ws.Range("E14").Value = "  -4.48%  "

# Row 15 - Wrapped liquid staked Ether 2.0
Set-TextValue "D15" "2.339.12"
$ws.Range("E15").Value = "  -1.94%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.33"
$ws.Range("E16").Value = "  -2.88%  "

# Row 17 - Wrapped Ether
Set-TextValue "D17" "2.045.81"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18 - Avalanche
Set-TextValue "D18" "17.50"
$ws.Range("E18").Value = "  +1.88%  "

# Row 19 - Wrapped BTC
Set-TextValue "D19" "36.196.61"
$ws.Range("E19").Value = "  -1.47%  "

# Row 20 - Litecoin
Set-TextValue "D20" "71.34"
$ws.Range("E20").Value = "  -1.74%  "

# Row 21 - Shiba Inu
$ws.Range("E21").Value = "  -2.73%  "

# Row 22 - Bitcoin Cash
Set-TextValue "D22" "235.96"
$ws.Range("E22").Value = "  -1.33%  "

# Row 23 - Uniswap
Set-TextValue "D23" "5.18"
$ws.Range("E23").Value = "  -4.48%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.15%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -2.21%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "2.28"
$ws.Range("E26").Value = "  +5.63%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.25"
$ws.Range("E27").Value = "  -5.77%  "

# Row 28 - Monero
Set-TextValue "D28" "163.49"
$ws.Range("E28").Value = "  -2.17%  "

# Row 29 - Ethereum Classic
Set-TextValue "D29" "19.84"
$ws.Range("E29").Value = "  -3.74%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  -1.70%  "

# Row 31 - ImmutableX
Set-TextValue "D31" "1.16"
$ws.Range("E31").Value = "  -1.27%  "

# Row 32 - Filecoin
Set-TextValue "D32" "4.95"
$ws.Range("E32").Value = "  -7.27%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0597"
$ws.Range("E33").Value = "  -1.34%  "

# Row 34 - Internet Computer (DFINITY)
$ws.Range("E34").Value = "  -6.31%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.0901"
$ws.Range("E35").Value = "  +8.56%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  -0.05%  "

# Row 37 - WEMIX Token
$ws.Range("E37").Value = "  -1.50%  "

# Row 38 - Lido DAO Token
Set-TextValue "D38" "2.20"
$ws.Range("E38").Value = "  -7.43%  "

# Row 39 - THORChain
Set-TextValue "D39" "5.03"
$ws.Range("E39").Value = "  +3.82%  "

# Row 40 - Trust Wallet Token
Set-TextValue "D40" "1.20"
$ws.Range("E40").Value = "  -5.03%  "

# Row 41 - Huobi Token
$ws.Range("E41").Value = "  +2.07%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  -2.61%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  -4.47%  "

# Row 44 / 45 - coins swapped: Aave <-> Cronos
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D44" "0.0903"
$ws.Range("E44").Value = "  -4.69%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "92.43"
$ws.Range("E45").Value = "  -4.13%  "

# Row 46 - Maker
Set-TextValue "D46" "1.401.02"
$ws.Range("E46").Value = "  +3.19%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  +2.11%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "15.53"
$ws.Range("E48").Value = "  -2.98%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  +2.10%  "

# Row 51 - MultiversX
Set-TextValue "D51" "45.72"
$ws.Range("E51").Value = "  +1.37%  "
